$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.347249388694763
$ws.Range("B1").Value = 2.334668159484863
$ws.Range("C1").Value = 2.835295677185059
$ws.Range("D1").Value = 3.268305063247681
$ws.Range("E1").Value = 1.970446467399597
